# Generate Report for Handback
# This script updates the localization-status workbook:
#  - Marks files as "Handed back: in sync with en-US" (was "Ready for handoff")
#    on the Overview sheet and on each language sheet's Status column.
#  - Stamps the Latest Handback DateTime for each language sheet.
#  - Adds "Latest Target File" (F) and "Latest Handback File" (G) hyperlink
#    columns' values for the two data rows on each language sheet.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Overview sheet: update the per-language status cells
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $handedBack
$overview.Range("C2").Value = $handedBack
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

# ---------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Status column
$zhcn.Range("C2").Value = $handedBack
$zhcn.Range("C3").Value = $handedBack

# Latest Handback DateTime
$zhcn.Range("H2").Value = "2016-03-12 02:44:05"
$zhcn.Range("H3").Value = "2016-03-12 02:44:05"

# Rebuild the hyperlinks in the desired order: existing A/B/D hyperlinks for
# row 2, the new F2/G2 hyperlinks, then the existing A/B/D hyperlinks for
# row 3, followed by the new F3/G3 hyperlinks.
$zhcn.Hyperlinks.Delete()

$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a51fe554778bd8409799c825fe154167d8a9e9cd/e2e/03824bdb-d6e6-4451-b5b5-471410216c3e.md", "", "", "03824bdb-d6e6-4451-b5b5-471410216c3e.md")
$zhcn.Hyperlinks.Add($zhcn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/a51fe554778bd8409799c825fe154167d8a9e9cd/e2e/03824bdb-d6e6-4451-b5b5-471410216c3e.md", "", "", ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68c5e549dc68b2d6e4cb2bb888353fbbb6c03221/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.zh-cn.xlf", "", "", "03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/a51fe554778bd8409799c825fe154167d8a9e9cd/e2e/03824bdb-d6e6-4451-b5b5-471410216c3e.md", "", "", "03824bdb-d6e6-4451-b5b5-471410216c3e.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68c5e549dc68b2d6e4cb2bb888353fbbb6c03221/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.zh-cn.xlf", "", "", "03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.zh-cn.xlf")

$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a51fe554778bd8409799c825fe154167d8a9e9cd/e2e/ffffb3aabc47-f91b-4db8-9836-c79007cce269.md", "", "", "ffffb3aabc47-f91b-4db8-9836-c79007cce269.md")
$zhcn.Hyperlinks.Add($zhcn.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/a51fe554778bd8409799c825fe154167d8a9e9cd/e2e/ffffb3aabc47-f91b-4db8-9836-c79007cce269.md", "", "", ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68c5e549dc68b2d6e4cb2bb888353fbbb6c03221/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.zh-cn.xlf", "", "", "03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/a51fe554778bd8409799c825fe154167d8a9e9cd/e2e/03824bdb-d6e6-4451-b5b5-471410216c3e.md", "", "", "03824bdb-d6e6-4451-b5b5-471410216c3e.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68c5e549dc68b2d6e4cb2bb888353fbbb6c03221/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.zh-cn.xlf", "", "", "03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.zh-cn.xlf")

# ---------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Status column
$dede.Range("C2").Value = $handedBack
$dede.Range("C3").Value = $handedBack

# Latest Handback DateTime (distinct stamp for de-de)
$dede.Range("H2").Value = "2016-03-12 02:44:10"
$dede.Range("H3").Value = "2016-03-12 02:44:10"

$dede.Hyperlinks.Delete()

$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a51fe554778bd8409799c825fe154167d8a9e9cd/e2e/03824bdb-d6e6-4451-b5b5-471410216c3e.md", "", "", "03824bdb-d6e6-4451-b5b5-471410216c3e.md")
$dede.Hyperlinks.Add($dede.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/a51fe554778bd8409799c825fe154167d8a9e9cd/e2e/03824bdb-d6e6-4451-b5b5-471410216c3e.md", "", "", ".md")
$dede.Hyperlinks.Add($dede.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bd8e4b8aa4ada1baad2f5f018623efae29f40a8d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.de-de.xlf", "", "", "03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/a51fe554778bd8409799c825fe154167d8a9e9cd/e2e/03824bdb-d6e6-4451-b5b5-471410216c3e.md", "", "", "03824bdb-d6e6-4451-b5b5-471410216c3e.md")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bd8e4b8aa4ada1baad2f5f018623efae29f40a8d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.de-de.xlf", "", "", "03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.de-de.xlf")

$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a51fe554778bd8409799c825fe154167d8a9e9cd/e2e/ffffb3aabc47-f91b-4db8-9836-c79007cce269.md", "", "", "ffffb3aabc47-f91b-4db8-9836-c79007cce269.md")
$dede.Hyperlinks.Add($dede.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/a51fe554778bd8409799c825fe154167d8a9e9cd/e2e/ffffb3aabc47-f91b-4db8-9836-c79007cce269.md", "", "", ".md")
$dede.Hyperlinks.Add($dede.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bd8e4b8aa4ada1baad2f5f018623efae29f40a8d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.de-de.xlf", "", "", "03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/a51fe554778bd8409799c825fe154167d8a9e9cd/e2e/03824bdb-d6e6-4451-b5b5-471410216c3e.md", "", "", "03824bdb-d6e6-4451-b5b5-471410216c3e.md")
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bd8e4b8aa4ada1baad2f5f018623efae29f40a8d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.de-de.xlf", "", "", "03824bdb-d6e6-4451-b5b5-471410216c3e.2239120a0a49db1a9e1ad4d040b293a8da2e7ba3.de-de.xlf")
